$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 533.3333
$ws.Range("I18").Value = 533.3333
$ws.Range("K18").Value = 533.3333
$ws.Range("M18").Value = -249.3333

$ws.Range("H33").Value = 10417325
$ws.Range("I33").Value = 15625492
$ws.Range("J33").Value = 990.125
$ws.Range("K33").Value = 15625492
$ws.Range("L33").Value = 990.125
$ws.Range("M33").Value = -15625263
$ws.Range("N33").Value = -1448.125

$ws.Range("H40").Value = 3666.6667
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""

$ws.Range("H80").Value = 3363.1667
$ws.Range("J80").Value = 6776.125
$ws.Range("L80").Value = 20328.375
$ws.Range("N80").Value = -22324.375

$ws.Range("H83").Value = 3363.1667
$ws.Range("J83").Value = 6776.125
$ws.Range("L83").Value = 60985.125
$ws.Range("N83").Value = -70969.125

$ws.Range("H137").Value = 10373.53
$ws.Range("I137").Value = 1697.3334
$ws.Range("J137").Value = 24072.79
$ws.Range("K137").Value = 5092.0002
$ws.Range("L137").Value = 72218.37
$ws.Range("M137").Value = -2542.0002
$ws.Range("N137").Value = -77318.37

$ws.Range("H138").Value = 4068.745
$ws.Range("I138").Value = 5200.9287
$ws.Range("J138").Value = 3640.3513
$ws.Range("K138").Value = 15602.7861
$ws.Range("L138").Value = 10921.0539
$ws.Range("M138").Value = -10462.7861
$ws.Range("N138").Value = -21201.0539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 126.25
$ws.Range("J5").Value = 135
$ws.Range("L5").Value = 135
$ws.Range("N5").Value = -359

$ws.Range("H53").Value = 7039
$ws.Range("I53").Value = 7039
$ws.Range("K53").Value = 7039
$ws.Range("M53").Value = -6357

$ws.Range("H122").Value = 953769.8
$ws.Range("I122").Value = 1334613.8
$ws.Range("J122").Value = 1660
$ws.Range("K122").Value = 4003841.4
$ws.Range("L122").Value = 4980
$ws.Range("M122").Value = -4001391.4
$ws.Range("N122").Value = -9880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 126.25
$ws.Range("J4").Value = 135
$ws.Range("L4").Value = 135
$ws.Range("N4").Value = -365

$ws.Range("H22").Value = 11459.4
$ws.Range("I22").Value = 12723.777
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 12723.777
$ws.Range("L22").Value = 80
$ws.Range("M22").Value = -12550.777
$ws.Range("N22").Value = -426

$ws.Range("H80").Value = 1531.8077
$ws.Range("I80").Value = 1048.8334
$ws.Range("J80").Value = 1945.7858
$ws.Range("K80").Value = 1048.8334
$ws.Range("L80").Value = 1945.7858
$ws.Range("M80").Value = -50.83339999999998
$ws.Range("N80").Value = -3941.7858

$ws.Range("H83").Value = 1531.8077
$ws.Range("I83").Value = 1048.8334
$ws.Range("J83").Value = 1945.7858
$ws.Range("K83").Value = 5244.166999999999
$ws.Range("L83").Value = 9728.929
$ws.Range("M83").Value = -252.1669999999995
$ws.Range("N83").Value = -19712.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 46610500
$ws.Range("I6").Value = 77683336
$ws.Range("J6").Value = 1250
$ws.Range("K6").Value = 77683336
$ws.Range("L6").Value = 1250
$ws.Range("M6").Value = -77683223
$ws.Range("N6").Value = -1476

$ws.Range("H16").Value = 7670.2666
$ws.Range("I16").Value = 9012.916999999999
$ws.Range("K16").Value = 9012.916999999999
$ws.Range("M16").Value = -8725.916999999999

$ws.Range("H22").Value = 1728
$ws.Range("I22").Value = 348.6
$ws.Range("K22").Value = 348.6
$ws.Range("M22").Value = 1.399999999999977

$ws.Range("H113").Value = 7670.2666
$ws.Range("I113").Value = 9012.916999999999
$ws.Range("K113").Value = 9012.916999999999
$ws.Range("M113").Value = -6842.916999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 491540.7
$ws.Range("I4").Value = 643159.7
$ws.Range("J4").Value = 91817.91
$ws.Range("K4").Value = 1929479.1
$ws.Range("L4").Value = 275453.73
$ws.Range("M4").Value = -1929367.1
$ws.Range("N4").Value = -275677.73

$ws.Range("H25").Value = 536.25
$ws.Range("I25").Value = 536.25
$ws.Range("K25").Value = 1608.75
$ws.Range("M25").Value = -1439.75

$ws.Range("H30").Value = 536.25
$ws.Range("I30").Value = 536.25
$ws.Range("K30").Value = 1608.75
$ws.Range("M30").Value = -1506.75

$ws.Range("H37").Value = 167998
$ws.Range("J37").Value = 167998
$ws.Range("L37").Value = 503994
$ws.Range("N37").Value = -504218

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 35.333332
$ws.Range("I2").Value = 35.333332
$ws.Range("K2").Value = 35.333332
$ws.Range("M2").Value = 77.666668

$ws.Range("H70").Value = 17142
$ws.Range("J70").Value = 15713.429
$ws.Range("L70").Value = 15713.429
$ws.Range("N70").Value = -16253.429

$ws.Range("H73").Value = 17142
$ws.Range("J73").Value = 15713.429
$ws.Range("L73").Value = 15713.429
$ws.Range("N73").Value = -17585.429

$ws.Range("H97").Value = 980.7143
$ws.Range("I97").Value = 748.5
$ws.Range("K97").Value = 748.5
$ws.Range("M97").Value = -252.5

$ws.Range("H107").Value = 710.6
$ws.Range("I107").Value = 763.25
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 763.25
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1156.75
$ws.Range("N107").Value = -4340

$ws.Range("H113").Value = 5999.6665
$ws.Range("I113").Value = 6999.5
$ws.Range("K113").Value = 6999.5
$ws.Range("M113").Value = -4829.5

$ws.Range("H135").Value = 240000
$ws.Range("J135").Value = 240000
$ws.Range("L135").Value = 240000
$ws.Range("N135").Value = -250140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2414.5
$ws.Range("I16").Value = 2310.7778
$ws.Range("K16").Value = 2310.7778
$ws.Range("M16").Value = -2140.7778

$ws.Range("H22").Value = 456.9
$ws.Range("I22").Value = 403.33334
$ws.Range("J22").Value = 537.25
$ws.Range("K22").Value = 403.33334
$ws.Range("L22").Value = 537.25
$ws.Range("M22").Value = -108.33334
$ws.Range("N22").Value = -1127.25

$ws.Range("H27").Value = 456.9
$ws.Range("I27").Value = 403.33334
$ws.Range("J27").Value = 537.25
$ws.Range("K27").Value = 403.33334
$ws.Range("L27").Value = 537.25
$ws.Range("M27").Value = -296.33334
$ws.Range("N27").Value = -751.25

$ws.Range("H68").Value = 5105913.5
$ws.Range("I68").Value = 2403
$ws.Range("K68").Value = 2403
$ws.Range("M68").Value = -1654

$ws.Range("H71").Value = 5105913.5
$ws.Range("I71").Value = 2403
$ws.Range("K71").Value = 12015
$ws.Range("M71").Value = -8271

$ws.Range("H122").Value = 20666
$ws.Range("I122").Value = 27999.5
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 83998.5
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = -81548.5
$ws.Range("N122").Value = -22897

$ws.Range("H136").Value = 2136219.5
$ws.Range("I136").Value = 27082.625
$ws.Range("J136").Value = 4546662
$ws.Range("K136").Value = 81247.875
$ws.Range("L136").Value = 13639986
$ws.Range("M136").Value = -78697.875
$ws.Range("N136").Value = -13645086

$ws.Range("H137").Value = 99285.71000000001
$ws.Range("J137").Value = 99285.71000000001
$ws.Range("L137").Value = 99285.71000000001
$ws.Range("N137").Value = -109485.71

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4353
$ws.Range("I62").Value = 3759.5
$ws.Range("J62").Value = 6331.3335
$ws.Range("K62").Value = 3759.5
$ws.Range("L62").Value = 6331.3335
$ws.Range("M62").Value = -3135.5
$ws.Range("N62").Value = -7579.3335

$ws.Range("H65").Value = 4353
$ws.Range("I65").Value = 3759.5
$ws.Range("J65").Value = 6331.3335
$ws.Range("K65").Value = 18797.5
$ws.Range("L65").Value = 31656.6675
$ws.Range("M65").Value = -15677.5
$ws.Range("N65").Value = -37896.6675

$ws.Range("H126").Value = 3593.5789
$ws.Range("I126").Value = 2861.3333
$ws.Range("J126").Value = 5390.909
$ws.Range("K126").Value = 8583.999899999999
$ws.Range("L126").Value = 16172.727
$ws.Range("M126").Value = -6113.999899999999
$ws.Range("N126").Value = -21112.727
